# hands-on-3+washup.pptx update: generalize "two PEs" wording to "multiple PEs"
# and add a second TO DO item about running on more than 2 PEs.

$p = $ppt.ActivePresentation

# --- Slide 2 ("Objective") ---------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange

$title = $tr2.Paragraphs(1, 1)
$tr2.Characters($title.Start, $title.Length).Text = "Perform Matrix-Vector Multiplication on multiple adjacent PEs"

$nCols = $tr2.Paragraphs(5, 1)
$tr2.Characters($nCols.Start, $nCols.Length).Text = "N columns will be split across the multiple PEs"

$xSplit = $tr2.Paragraphs(7, 1)
$tr2.Characters($xSplit.Start, $xSplit.Length).Text = "x will be split across multiple PEs"

# --- Slide 3 ("To do:") --------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange

# Paragraph 2 is "TO DO 1: Fill in mpi_x details for reduction".
# Add a blank line followed by a new "TO DO 2" line straight after it,
# ahead of the pre-existing trailing blank paragraph.
$toDo1 = $tr3.Paragraphs(2, 1)
$toDo1.InsertAfter("`r`r") | Out-Null

$toDo2 = $tr3.Paragraphs(4, 1)
$toDo2.Text = "TO DO 2: What needs to be changed to run on more than 2, say 4 PEs?"
